$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "42.321.33"
$ws.Range("E2").Value = "  -3.41%  "
Set-TextValue $ws "D3" "2.245.83"
$ws.Range("E3").Value = "  -4.12%  "
$ws.Range("E4").Value = "  -0.19%  "
Set-TextValue $ws "D5" "233.59"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("E6").Value = "  -5.12%  "
Set-TextValue $ws "D7" "69.69"
$ws.Range("E7").Value = "  -3.24%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue $ws "D9" "0.561"
$ws.Range("E9").Value = "  -4.83%  "
Set-TextValue $ws "D10" "0.0991"
$ws.Range("E10").Value = "  -0.26%  "
Set-TextValue $ws "D11" "58.33"
$ws.Range("E11").Value = "  -0.19%  "
Set-TextValue $ws "D12" "35.98"
$ws.Range("E12").Value = "  +11.21%  "
$ws.Range("E13").Value = "  -2.06%  "
Set-TextValue $ws "D14" "6.76"
$ws.Range("E14").Value = "  -5.40%  "
Set-TextValue $ws "D15" "2.583.31"
$ws.Range("E15").Value = "  -4.12%  "
Set-TextValue $ws "D16" "15.05"
$ws.Range("E16").Value = "  -7.07%  "
Set-TextValue $ws "D17" "0.860"
$ws.Range("E17").Value = "  -4.34%  "
Set-TextValue $ws "D18" "2.247.17"
$ws.Range("E18").Value = "  -4.29%  "
Set-TextValue $ws "D19" "42.231.25"
$ws.Range("E19").Value = "  -3.38%  "
Set-TextValue $ws "D20" "0.0₃0975"
$ws.Range("E20").Value = "  -3.66%  "
Set-TextValue $ws "D21" "6.27"
$ws.Range("E21").Value = "  -6.06%  "
$ws.Range("E22").Value = "  -5.81%  "
Set-TextValue $ws "D23" "236.50"
$ws.Range("E23").Value = "  -6.58%  "
Set-TextValue $ws "D24" "1.98"
$ws.Range("E24").Value = "  +3.89%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -4.67%  "
Set-TextValue $ws "D28" "10.03"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("E29").Value = "  -4.01%  "
Set-TextValue $ws "D30" "169.25"
$ws.Range("E30").Value = "  -3.74%  "
Set-TextValue $ws "D31" "20.57"
$ws.Range("E31").Value = "  -7.42%  "
$ws.Range("E32").Value = "  -4.73%  "
Set-TextValue $ws "D33" "0.127"
$ws.Range("E33").Value = "  -6.08%  "
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("E35").Value = "  -0.11%  "
Set-TextValue $ws "D36" "4.72"
$ws.Range("E36").Value = "  -6.96%  "
$ws.Range("E37").Value = "  -2.73%  "
Set-TextValue $ws "D38" "21.53"
$ws.Range("E38").Value = "  +14.21%  "
$ws.Range("E39").Value = "  -4.67%  "
Set-TextValue $ws "D40" "0.0268"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  -6.53%  "
Set-TextValue $ws "D42" "66.21"
$ws.Range("E42").Value = "  +1.84%  "
Set-TextValue $ws "D43" "4.91"
$ws.Range("E43").Value = "  -7.01%  "
Set-TextValue $ws "D44" "8.95"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("E45").Value = "  -4.39%  "
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws "D46" "1.00"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D47" "0.189"
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue $ws "D48" "4.43"
$ws.Range("E48").Value = "  +10.21%  "
$ws.Range("B49").Value = "BitTorrent-New"
$ws.Range("C49").Value = "https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt"
Set-TextValue $ws "D49" "0.0₃0151"
$ws.Range("E49").Value = "  +14.91%  "
$ws.Range("E50").Value = "  -3.63%  "
Set-TextValue $ws "D51" "2.34"
$ws.Range("E51").Value = "  -4.33%  "
